$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B=-0.09347130405523991; C=1.163400873526787;  D=8.140377539043413; E=2.853134686453378; F=2.883113556558969; G=46},
    @{Row=3;  B=0.02024706560588739;  C=1.092969939287157;  D=5.366131528809891; E=2.31649121060493;  F=2.342577598071925; G=45},
    @{Row=4;  B=-0.03266126848104425; C=0.9435118373712057; D=4.478804720633087; E=2.116318671805616; F=2.14053063588333;  G=44},
    @{Row=5;  B=0.07925189262698597; C=1.009367654643358;   D=4.81238600652716;  E=2.193715115170418; F=2.218228171063641; G=43},
    @{Row=6;  B=0.02771576539098436; C=0.9479935622673827;  D=4.345380084005096; E=2.084557527151768; F=2.109639321010291; G=42},
    @{Row=7;  B=0.1037329126159212;  C=0.9868912652243453;  D=4.448777381052547; E=2.109212502583025; F=2.132830815017939; G=41},
    @{Row=8;  B=0.06410150100714884; C=0.9537641320944577;  D=4.474893186289131; E=2.115394333520143; F=2.141359305616774; G=40},
    @{Row=9;  B=0.1191595648948872;  C=1.014371866898358;   D=4.568013714778091; E=2.137291209633842; F=2.16186307201873;  G=39},
    @{Row=10; B=0.08350116669570022; C=0.9748485792014643;  D=4.600399335576155; E=2.14485415251857;  F=2.171997597480076; G=38},
    @{Row=11; B=0.116349882878518;   C=1.0065703756427;     D=4.725732831804032; E=2.173875072722449; F=2.200702128189578; G=37}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.B
    $ws.Range("C$r").Value = $item.C
    $ws.Range("D$r").Value = $item.D
    $ws.Range("E$r").Value = $item.E
    $ws.Range("F$r").Value = $item.F
    $ws.Range("G$r").Value = $item.G
}
